# Auto-generated edit script: updates market-price derived columns (H:N)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR leve-profit sheets.
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 300
$ws.Range("I29").Value = 300
$ws.Range("K29").Value = 900
$ws.Range("M29").Value = -619
$ws.Range("H58").Value = 4462.476
$ws.Range("I58").Value = 104.4
$ws.Range("J58").Value = 8424.362999999999
$ws.Range("K58").Value = 313.2
$ws.Range("L58").Value = 25273.089
$ws.Range("M58").Value = -163.2
$ws.Range("N58").Value = -25573.089
$ws.Range("H92").Value = 22222828
$ws.Range("I92").Value = 27778034
$ws.Range("J92").Value = 2005
$ws.Range("K92").Value = 27778034
$ws.Range("L92").Value = 2005
$ws.Range("M92").Value = -27776786
$ws.Range("N92").Value = -4501
$ws.Range("H98").Value = 433352.03
$ws.Range("I98").Value = 624985.6
$ws.Range("J98").Value = 2176.5
$ws.Range("K98").Value = 624985.6
$ws.Range("L98").Value = 2176.5
$ws.Range("M98").Value = -623487.6
$ws.Range("N98").Value = -5172.5
$ws.Range("H106").Value = 18521370
$ws.Range("J106").Value = 3000
$ws.Range("L106").Value = 3000
$ws.Range("N106").Value = -4262
$ws.Range("H113").Value = 4153.1333
$ws.Range("I113").Value = 3613.8572
$ws.Range("J113").Value = 4625
$ws.Range("K113").Value = 3613.8572
$ws.Range("L113").Value = 4625
$ws.Range("M113").Value = -359.8571999999999
$ws.Range("N113").Value = -11133
$ws.Range("H116").Value = 2009.75
$ws.Range("J116").Value = 1377
$ws.Range("L116").Value = 1377
$ws.Range("N116").Value = -8261
$ws.Range("H122").Value = 433352.03
$ws.Range("I122").Value = 624985.6
$ws.Range("J122").Value = 2176.5
$ws.Range("K122").Value = 1874956.8
$ws.Range("L122").Value = 6529.5
$ws.Range("M122").Value = -1872506.8
$ws.Range("N122").Value = -11429.5
$ws.Range("H132").Value = 15381.069
$ws.Range("I132").Value = 17974.377
$ws.Range("J132").Value = 1000
$ws.Range("K132").Value = 53923.131
$ws.Range("L132").Value = 3000
$ws.Range("M132").Value = -51393.131
$ws.Range("N132").Value = -8060
$ws.Range("H138").Value = 5894142
$ws.Range("I138").Value = 2103350.5
$ws.Range("J138").Value = 7939966
$ws.Range("K138").Value = 6310051.5
$ws.Range("L138").Value = 23819898
$ws.Range("M138").Value = -6304911.5
$ws.Range("N138").Value = -23830178

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 24414.129
$ws.Range("I32").Value = 6648.4385
$ws.Range("J32").Value = 102309.84
$ws.Range("K32").Value = 6648.4385
$ws.Range("L32").Value = 102309.84
$ws.Range("M32").Value = -6361.4385
$ws.Range("N32").Value = -102883.84
$ws.Range("H45").Value = 1357.7778
$ws.Range("I45").Value = 1444.4286
$ws.Range("J45").Value = 1054.5
$ws.Range("K45").Value = 1444.4286
$ws.Range("L45").Value = 1054.5
$ws.Range("M45").Value = -1067.4286
$ws.Range("N45").Value = -1808.5
$ws.Range("H122").Value = 3676.158
$ws.Range("I122").Value = 3597.0557
$ws.Range("J122").Value = 5100
$ws.Range("K122").Value = 10791.1671
$ws.Range("L122").Value = 15300
$ws.Range("M122").Value = -8341.167099999999
$ws.Range("N122").Value = -20200

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 505755
$ws.Range("I3").Value = 505755
$ws.Range("K3").Value = 505755
$ws.Range("M3").Value = -505641
$ws.Range("H20").Value = 1892.4286
$ws.Range("J20").Value = 2100
$ws.Range("L20").Value = 2100
$ws.Range("N20").Value = -2594
$ws.Range("H105").Value = 2881.932
$ws.Range("I105").Value = 2722.3872
$ws.Range("J105").Value = 3262.3845
$ws.Range("K105").Value = 2722.3872
$ws.Range("L105").Value = 3262.3845
$ws.Range("M105").Value = -975.3872000000001
$ws.Range("N105").Value = -6756.3845
$ws.Range("H132").Value = 45640
$ws.Range("J132").Value = 45640
$ws.Range("L132").Value = 45640
$ws.Range("N132").Value = -55760

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").Value = ""
$ws.Range("H122").Value = 1262.7
$ws.Range("I122").Value = 1234.421
$ws.Range("J122").Value = 1800
$ws.Range("K122").Value = 3703.263
$ws.Range("L122").Value = 5400
$ws.Range("M122").Value = -1253.263
$ws.Range("N122").Value = -10300
$ws.Range("H132").Value = 3335110.2
$ws.Range("I132").Value = 4387164
$ws.Range("J132").Value = 3606.5
$ws.Range("K132").Value = 13161492
$ws.Range("L132").Value = 10819.5
$ws.Range("M132").Value = -13158962
$ws.Range("N132").Value = -15879.5

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 4488.8887
$ws.Range("H69").Value = 5857.143
$ws.Range("I69").Value = 1000
$ws.Range("J69").Value = 6666.6665
$ws.Range("K69").Value = 3000
$ws.Range("L69").Value = 19999.9995
$ws.Range("M69").Value = -2189
$ws.Range("N69").Value = -21621.9995
$ws.Range("H72").Value = 5857.143
$ws.Range("I72").Value = 1000
$ws.Range("J72").Value = 6666.6665
$ws.Range("K72").Value = 9000
$ws.Range("L72").Value = 59999.9985
$ws.Range("M72").Value = -4944
$ws.Range("N72").Value = -68111.9985
$ws.Range("H103").Value = 1975.1666
$ws.Range("I103").Value = 1089.1111
$ws.Range("J103").Value = 4633.3335
$ws.Range("K103").Value = 3267.3333
$ws.Range("L103").Value = 13900.0005
$ws.Range("M103").Value = -2388.3333
$ws.Range("N103").Value = -15658.0005
$ws.Range("H127").Value = 1476.2778
$ws.Range("J127").Value = 1527.8235
$ws.Range("L127").Value = 4583.470499999999
$ws.Range("N127").Value = -14503.4705
$ws.Range("H131").Value = 13891024
$ws.Range("J131").Value = 15875326
$ws.Range("L131").Value = 47625978
$ws.Range("N131").Value = -47636058

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 2716000
$ws.Range("I11").Value = 6333333.5
$ws.Range("K11").Value = 6333333.5
$ws.Range("M11").Value = -6333194.5
$ws.Range("H122").Value = 2207.6296
$ws.Range("I122").Value = 2295.75
$ws.Range("J122").Value = 1502.6666
$ws.Range("K122").Value = 6887.25
$ws.Range("L122").Value = 4507.9998
$ws.Range("M122").Value = -4437.25
$ws.Range("N122").Value = -9407.9998

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3131.76
$ws.Range("I7").Value = 2454.889
$ws.Range("J7").Value = 3512.5
$ws.Range("K7").Value = 2454.889
$ws.Range("L7").Value = 3512.5
$ws.Range("M7").Value = -2342.889
$ws.Range("N7").Value = -3736.5
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("M26").Value = ""
$ws.Range("H40").Value = 2862.0952
$ws.Range("I40").Value = 1825.5
$ws.Range("J40").Value = 3500
$ws.Range("K40").Value = 1825.5
$ws.Range("L40").Value = 3500
$ws.Range("M40").Value = -1689.5
$ws.Range("N40").Value = -3772
$ws.Range("H126").Value = 3131.76
$ws.Range("I126").Value = 2454.889
$ws.Range("J126").Value = 3512.5
$ws.Range("K126").Value = 7364.667
$ws.Range("L126").Value = 10537.5
$ws.Range("M126").Value = -4894.667
$ws.Range("N126").Value = -15477.5
$ws.Range("H132").Value = 5182.7646
$ws.Range("I132").Value = 3727.5715
$ws.Range("J132").Value = 6201.4
$ws.Range("K132").Value = 11182.7145
$ws.Range("L132").Value = 18604.2
$ws.Range("M132").Value = -8652.7145
$ws.Range("N132").Value = -23664.2

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1893.3334
$ws.Range("I122").Value = 1707.6923
$ws.Range("K122").Value = 5123.0769
$ws.Range("M122").Value = -2673.0769
$ws.Range("H126").Value = 91720.17999999999
$ws.Range("I126").Value = 125515.375
$ws.Range("J126").Value = 1599.6666
$ws.Range("K126").Value = 376546.125
$ws.Range("L126").Value = 4798.9998
$ws.Range("M126").Value = -374076.125
$ws.Range("N126").Value = -9738.9998
$ws.Range("H132").Value = 5018.483
$ws.Range("I132").Value = 5754.1333
$ws.Range("K132").Value = 17262.3999
$ws.Range("M132").Value = -14732.3999
$ws.Range("H136").Value = 2371.0637
$ws.Range("I136").Value = 640.86487
$ws.Range("J136").Value = 8772.799999999999
$ws.Range("K136").Value = 1922.59461
$ws.Range("L136").Value = 26318.4
$ws.Range("M136").Value = 627.4053899999999
$ws.Range("N136").Value = -31418.4

